# Update the "想去人数" (F column) figures on both the "展览" and "全部类型"
# worksheets to match the newly scraped counts.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 1291
    $ws.Range("F3").Value = 1681
    $ws.Range("F5").Value = 6228
    $ws.Range("F6").Value = 63
}
